$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-13 (columns A-T), reflecting the updated
# NATMI Cx3cl1-Itgav LR-pair results (re-run per Dr Hou advice), including
# the newly introduced "sCs" sending-cluster rows 10-13.
$data = @(
  @("ECs", "Cx3cl1", "Itgav", "ECs", 2, 0.6666666666666666, 3.744854333333334, 11.234563, 0.2406267676902895, 0.2406267676902895, 3, 1, 15.35884066666667, 46.076522, 0.1012042817263867, 0.1012042817263867, 57.5166210255429, 517.649589229886, 0.02435245918823785, 0.02435245918823785),
  @("ECs", "Cx3cl1", "Itgav", "FAPs", 2, 0.6666666666666666, 3.744854333333334, 11.234563, 0.2406267676902895, 0.2406267676902895, 3, 1, 50.59256466666667, 151.777694, 0.3333704853712116, 0.3333704853712116, 189.4617850264136, 1705.156065237722, 0.08021786233821757, 0.08021786233821757),
  @("ECs", "Cx3cl1", "Itgav", "M2", 2, 0.6666666666666666, 3.744854333333334, 11.234563, 0.2406267676902895, 0.2406267676902895, 3, 1, 60.37715666666667, 181.13147, 0.397844271305776, 0.397844271305776, 226.1036567775123, 2034.93291099761, 0.09573198104840745, 0.09573198104840745),
  @("ECs", "Cx3cl1", "Itgav", "sCs", 2, 0.6666666666666666, 3.744854333333334, 11.234563, 0.2406267676902895, 0.2406267676902895, 3, 1, 25.43221733333333, 76.296652, 0.1675809615966257, 0.1675809615966258, 95.23994928700844, 857.159543583076, 0.04032446511542658, 0.04032446511542659),
  @("FAPs", "Cx3cl1", "Itgav", "ECs", 3, 1, 10.317205, 30.951615, 0.662935182458301, 0.662935182458301, 3, 1, 15.35884066666667, 46.076522, 0.1012042817263867, 0.1012042817263867, 158.4603077203367, 1426.14276948303, 0.06709187897184343, 0.06709187897184343),
  @("FAPs", "Cx3cl1", "Itgav", "FAPs", 3, 1, 10.317205, 30.951615, 0.662935182458301, 0.662935182458301, 3, 1, 50.59256466666667, 151.777694, 0.3333704853712116, 0.3333704853712116, 521.9738611417566, 4697.764750275809, 0.2210030235457765, 0.2210030235457765),
  @("FAPs", "Cx3cl1", "Itgav", "M2", 3, 1, 10.317205, 30.951615, 0.662935182458301, 0.662935182458301, 3, 1, 60.37715666666667, 181.13147, 0.397844271305776, 0.397844271305776, 622.9235026471167, 5606.31152382405, 0.2637449645880844, 0.2637449645880844),
  @("FAPs", "Cx3cl1", "Itgav", "sCs", 3, 1, 10.317205, 30.951615, 0.662935182458301, 0.662935182458301, 3, 1, 25.43221733333333, 76.296652, 0.1675809615966257, 0.1675809615966258, 262.3893998325533, 2361.50459849298, 0.1110953153525966, 0.1110953153525966),
  @("sCs", "Cx3cl1", "Itgav", "ECs", 3, 1, 1.500857333333333, 4.502572, 0.0964380498514096, 0.0964380498514096, 3, 1, 15.35884066666667, 46.076522, 0.1012042817263867, 0.1012042817263867, 23.05142864606488, 207.462857814584, 0.009759943566305378, 0.009759943566305378),
  @("sCs", "Cx3cl1", "Itgav", "FAPs", 3, 1, 1.500857333333333, 4.502572, 0.0964380498514096, 0.0964380498514096, 3, 1, 50.59256466666667, 151.777694, 0.3333704853712116, 0.3333704853712116, 75.93222169210756, 683.3899952289679, 0.03214959948721751, 0.03214959948721751),
  @("sCs", "Cx3cl1", "Itgav", "M2", 3, 1, 1.500857333333333, 4.502572, 0.0964380498514096, 0.0964380498514096, 3, 1, 60.37715666666667, 181.13147, 0.397844271305776, 0.397844271305776, 90.61749834898222, 815.55748514084, 0.03836732566928415, 0.03836732566928415),
  @("sCs", "Cx3cl1", "Itgav", "sCs", 3, 1, 1.500857333333333, 4.502572, 0.0964380498514096, 0.0964380498514096, 3, 1, 25.43221733333333, 76.296652, 0.1675809615966257, 0.1675809615966258, 38.17012988766044, 343.5311689889439, 0.01616118112860255, 0.01616118112860255),
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $rowNum = $i + 2
  $rowVals = $data[$i]
  for ($j = 0; $j -lt $rowVals.Length; $j++) {
    $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
  }
}
